$d = $word.ActiveDocument

$replacements = @(
    @('810÷4=', '688÷8='),
    @('869÷9=', '268÷8='),
    @('152÷4=', '307÷4='),
    @('584÷7=', '628÷4='),
    @('480÷9=', '646÷8='),
    @('710÷5=', '786÷2='),
    @('168÷6=', '960÷7='),
    @('329÷3=', '982÷9='),
    @('304÷4=', '918÷4='),
    @('280÷2=', '219÷9='),
    @('205÷4=', '544÷8='),
    @('705÷7=', '272÷9='),
    @('223÷9=', '222÷5='),
    @('937÷3=', '674÷4='),
    @('188÷4=', '425÷5='),
    @('103÷4=', '568÷8='),
    @('315÷3=', '617÷8='),
    @('227÷8=', '854÷9='),
    @('451÷3=', '249÷6='),
    @('318÷4=', '994÷8='),
    @('678÷2=', '560÷7='),
    @('733÷8=', '635÷8='),
    @('842÷9=', '209÷3='),
    @('268÷6=', '520÷6='),
    @('330÷5=', '101÷3='),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
